# Fixes some review findings in the SORMAS User Rights matrix:
# a number of cell pairs on the "User Rights" sheet had their
# Yes/No (green/red) values accidentally swapped; this script
# swaps each such pair back (value + fill/style together) using
# a scratch cell so that both the shared-string value and the
# cell's style index travel together, exactly as Excel's
# Cut/Copy-Paste would do it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Rights")

# Scratch cell used as temporary holding area while swapping two
# cells' full contents (value + formatting).
$scratch = $ws.Range("AZ1")

$swapPairs = @(
    "O3:P3",
    "C4:D4",
    "O4:P4",
    "O5:P5",
    "I7:J7",
    "O9:P9",
    "O11:P11",
    "I14:J14",
    "O14:P14",
    "I15:J15",
    "O19:P19",
    "O20:P20",
    "O21:P21",
    "O24:P24",
    "O25:P25",
    "O29:P29",
    "O35:P35",
    "O45:P45",
    "O48:P48",
    "O52:P52",
    "O60:P60",
    "I61:J61",
    "O61:P61",
    "I64:J64",
    "O64:P64",
    "O71:P71",
    "I74:J74",
    "O74:P74",
    "O80:P80",
    "O83:P83",
    "I95:J95",
    "O95:P95",
    "I96:J96",
    "O96:P96",
    "I106:J106",
    "I107:J107",
    "I110:J110",
    "O110:P110",
    "I113:J113",
    "O114:P114",
    "O127:P127",
    "I134:J134",
    "C136:D136",
    "E136:F136",
    "O136:P136",
    "C137:D137",
    "E137:F137",
    "C138:D138",
    "E138:F138",
    "O138:P138",
    "C139:D139",
    "E139:F139",
    "I140:J140",
    "I144:J144",
    "I145:J145",
    "I148:J148",
    "C149:D149",
    "I150:J150",
    "O153:P153",
    "O154:P154",
    "O155:P155",
    "O156:P156",
    "I160:J160",
    "I162:J162",
    "I163:J163",
    "C164:D164",
    "E165:F165"
)

foreach ($pair in $swapPairs) {
    $refs = $pair.Split(":")
    $cellA = $ws.Range($refs[0])
    $cellB = $ws.Range($refs[1])

    $cellA.Copy($scratch)
    $cellB.Copy($cellA)
    $scratch.Copy($cellB)
}

$scratch.Clear()
